# Apply META DATA batch list updates to Sheet1
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update fBatchId (L), fClassId (M), fDeptId (O) for student rows 2-11
for ($r = 2; $r -le 11; $r++) {
    $ws.Range("L$r").Value = 8
    $ws.Range("M$r").Value = 4
    $ws.Range("O$r").Value = 7
}

# Update the active selection to L11 to match the saved view state
$ws.Activate()
$ws.Range("L11").Select()
